$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.110.65"
$ws.Range("E2").Value = "  -0.24%  "
$ws.Range("D3").Value = "2.562.84"
$ws.Range("E3").Value = "  +0.45%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "585.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.99%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "148.00"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.597"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.61%  "
$ws.Range("E9").Value = "  +2.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "5.65"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.72%  "
$ws.Range("E11").Value = "  -0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.357"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.01%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "27.43"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.61%  "
$ws.Range("D14").Value = "3.022.86"
$ws.Range("E14").Value = "  +0.57%  "
$ws.Range("D15").Value = "63.038.64"
$ws.Range("E15").Value = "  -0.24%  "
$ws.Range("E16").Value = "  +2.17%  "
$ws.Range("D17").Value = "2.565.61"
$ws.Range("E17").Value = "  +0.83%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "11.36"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.10%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "344.05"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.45%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.88"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.32%  "
$ws.Range("E22").Value = "  +0.06%  "
$ws.Range("E23").Value = "  -3.63%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "66.55"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "2.678.52"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").Value = "  +0.22%  "
$ws.Range("E27").Value = "  -0.69%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.23"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.56%  "
$ws.Range("B29").Value = "SuiNetwork"
$ws.Range("C29").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.49"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.28%  "
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "8.43"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.07%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.64%  "
$ws.Range("D33").Value = "0.0₃0828"
$ws.Range("E33").Value = "  +0.36%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "462.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +11.95%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "175.94"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +2.63%  "
$ws.Range("E37").Value = "  +1.38%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "19.18"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.94%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.56"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.27%  "
$ws.Range("E40").Value = "  +0.03%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.56%  "
$ws.Range("E42").Value = "  +0.01%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "150.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.32%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.83"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.99%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "20.99"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.20%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0548"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.43%  "
$ws.Range("E47").Value = "  +1.13%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0976"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.29%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0241"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.45%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.75"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.12%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "11.39"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.43%  "
